{"js": "// Office.js (Word JavaScript API) edit script\n// Applies the diff: updates paragraphs 0-6 text, removes paragraphs 7-9,\n// and updates the final URL paragraph.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Paragraph 0: title line + subtitle line (separated by a manual line break).\nparagraphs.items[0].insertText(\n  \"\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 30.09.24: \u26a1\ufe0f\ud83d\ude80\" + \"\\u000b\" + \"SCHRODINGER\u2019S MEMORY: LARGE LANGUAGE MODELS\",\n  Word.InsertLocation.replace\n);\n\n// Paragraph 1: intro blurb.\nparagraphs.items[1].insertText(\"\u05d1\u05d9\u05d5\u05dd \u05d4\u05e1\u05d5\u05e2\u05e8 \u05d4\u05d6\u05d4 (\u05dc\u05de\u05e8\u05d5\u05ea \u05e9\u05d4\u05e1\u05e7\u05d9\u05e8\u05d4 \u05e9\u05d9\u05d9\u05db\u05ea \u05e4\u05d5\u05e8\u05de\u05dc\u05d9\u05ea \u05dc\u05d0\u05ea\u05de\u05d5\u05dc - \u05d0\u05e9\u05dc\u05d9\u05dd \u05d0\u05ea \u05d4\u05e4\u05e2\u05e8 \u05d1\u05d9\u05de\u05d9\u05dd \u05d4\u05e7\u05e8\u05d5\u05d1\u05d9\u05dd) \u05e0\u05e1\u05e7\u05d5\u05e8 \u05de\u05d0\u05de\u05e8 \u05d3\u05d9 \u05e7\u05dc\u05d9\u05dc \u05e2\u05dd \u05e9\u05dd \u05de\u05d0\u05d5\u05d3 \u05dc\u05d0 \u05e7\u05dc\u05d9\u05dc. \u05db\u05d9 \u05d0\u05d9\u05df \u05d3\u05d1\u05e8 \u05e7\u05dc\u05d9\u05dc \u05e9\u05db\u05d5\u05dc\u05dc \u05d1\u05ea\u05d5\u05db\u05d5 \u05d0\u05ea \u05e9\u05de\u05d5 \u05e9\u05dc \u05e9\u05e8\u05d3\u05d9\u05e0\u05d2\u05e8 - \u05e1\u05e4\u05e7 \u05d0\u05dd \u05d4\u05e6\u05dc\u05d7\u05ea\u05d9 \u05dc\u05d4\u05d1\u05d9\u05df \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d8\u05d5\u05d1\u05d4 \u05de\u05e1\u05e4\u05d9\u05e7 \u05d0\u05ea \u05d4\u05de\u05e9\u05d5\u05d5\u05d0\u05d4 \u05e9\u05dc \u05e9\u05e8\u05d3\u05d9\u05e0\u05d2\u05e8 \u05e2\u05d5\u05d3 \u05d1\u05e7\u05d5\u05e8\u05e1 \u05e4\u05d9\u05d6\u05d9\u05e7\u05d4 3 \u05d1\u05d0\u05d5\u05e0\u05d9\u05d1\u05e8\u05e1\u05d9\u05d8\u05d4 \u05d1\u05de\u05d5\u05e1\u05e7\u05d1\u05d4 \u05dc\u05e4\u05e0\u05d9 \u05e2\u05e9\u05e8\u05d5\u05ea \u05e9\u05e0\u05d9\u05dd. \u05d2\u05dd \u05e1\u05d9\u05e4\u05d5\u05e8\u05d5 \u05e9\u05dc \u05d7\u05ea\u05d5\u05dc \u05e9\u05e8\u05d3\u05d9\u05e0\u05d2\u05e8 \u05dc\u05d0 \u05d4\u05ea\u05d1\u05d4\u05e8 \u05e2\u05d3 \u05e2\u05db\u05e9\u05d9\u05d5.\", Word.InsertLocation.replace);\n\n// Paragraph 2: paper topic summary.\nparagraphs.items[2].insertText(\"\u05d0\u05d5\u05e7\u05d9\u05d9, \u05e1\u05d9\u05d9\u05de\u05e0\u05d5 \u05e2\u05dd \u05d4\u05e6\u05d7\u05d5\u05e7\u05d9\u05dd. \u05d4\u05de\u05d0\u05de\u05e8 \u05d7\u05d5\u05e7\u05e8 (\u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05ea) \u05e0\u05d5\u05e9\u05d0 \u05d3\u05d9 \u05e8\u05e6\u05d9\u05e0\u05d9 \u05d5\u05d4\u05d5\u05d0 \u05d4\u05d6\u05db\u05e8\u05d5\u05df \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4. \u05db\u05e9\u05d0\u05e0\u05d7\u05e0\u05d5 \u05e9\u05d5\u05d0\u05dc\u05d9\u05dd LLM \u05de\u05d4 \u05e2\u05d9\u05e8 \u05d4\u05d1\u05d9\u05e8\u05d4 \u05e9\u05dc \u05e9\u05d1\u05d3\u05d9\u05d4, \u05d0\u05d9\u05da \u05d4\u05d5\u05d0 \u05d9\u05d5\u05d3\u05e2 \u05e9\u05d6\u05d4 \u05e1\u05d8\u05d5\u05e7\u05d4\u05d5\u05dc\u05dd. \u05d4\u05de\u05d0\u05de\u05e8 \u05d8\u05d5\u05e2\u05df \u05db\u05d9 \u05d6\u05d9\u05db\u05e8\u05d5\u05df LLM \u05e4\u05d5\u05e2\u05dc \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05ea\u05d0\u05de\u05d4 \u05d3\u05d9\u05e0\u05de\u05d9\u05ea \u05e9\u05dc \u05e4\u05dc\u05d8\u05d9\u05dd \u05dc\u05e7\u05dc\u05d8\u05d9\u05dd. \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05de\u05d5\u05d3\u05dc \u05f4\u05d1\u05d5\u05d7\u05e8\u05f4 \u05d0\u05d9\u05da \u05dc\u05e9\u05dc\u05d5\u05e3 \u05d0\u05ea \u05d4\u05de\u05d9\u05d3\u05e2 \u05de\u05d4\u05d6\u05d9\u05db\u05e8\u05d5\u05df \u05d5\u05d1\u05d5\u05e0\u05d4 \u05d0\u05d5\u05ea\u05d5 \u05e2\u05dc \u05e1\u05de\u05da \u05d4\u05e7\u05dc\u05d8. \", Word.InsertLocation.replace);\n\n// Paragraph 3: authors' explanation (transformer architecture).\nparagraphs.items[3].insertText(\"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e1\u05d1\u05d9\u05e8\u05d9\u05dd \u05d0\u05ea \u05d0\u05d9\u05da \u05e4\u05d5\u05e2\u05dc \u05d4\u05d6\u05d9\u05db\u05e8\u05d5\u05df \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05e0\u05d9\u05ea\u05d5\u05d7 \u05e9\u05dc \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05ea \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd. \u05de\u05e0\u05d2\u05e0\u05d5\u05df \u05d4-attention (\u05db\u05dc\u05d5\u05de\u05e8 \u05de\u05e7\u05d3\u05de\u05d9 \u05d4-attention \u05e9\u05dc\u05d5) \u05dc\u05de\u05e2\u05e9\u05d4 \u05de\u05d0\u05e4\u05e9\u05e8\u05d9\u05dd \u05dc\u05de\u05d5\u05d3\u05dc \u05dc\u05d1\u05e0\u05d5\u05ea \u05d0\u05ea \u05d4\u05e4\u05dc\u05d8 \u05db\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05d3\u05d9\u05e0\u05de\u05d9\u05ea \u05e9\u05dc \u05d4\u05e7\u05dc\u05d8 (\u05db\u05dc\u05d5\u05de\u05e8 \u05dc\u05d0 \u05e7\u05d1\u05d5\u05e2\u05d4 \u05db\u05de\u05d5 \u05d1-MLP \u05d0\u05d5 ConvNets). \", Word.InsertLocation.replace);\n\n// Paragraph 4: Universal Approximation Theorem paragraph.\nparagraphs.items[4].insertText(\"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e9\u05ea\u05de\u05e9\u05d9\u05dd \u05d1- Universal Approximation Theorem \u05d0\u05d5 UAT \u05db\u05d3\u05d9 \u05dc\u05d4\u05e1\u05d1\u05d9\u05e8 \u05d0\u05ea \u05d4\u05d9\u05db\u05d5\u05dc\u05ea \u05e9\u05dc \u05e9\u05dc\u05d9\u05e4\u05ea \u05de\u05d9\u05d3\u05e2 \u05e9\u05e0\u05dc\u05de\u05d3 \u05d1\u05de\u05d4\u05dc\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05e2\u05dc \u05d1\u05e1\u05d9\u05e1 \u05ea\u05d5\u05db\u05df \u05e9\u05dc \u05d4\u05e7\u05dc\u05d8. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d8\u05d5\u05e2\u05e0\u05d9\u05dd \u05db\u05d9 \u05e0\u05d9\u05ea\u05df \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05e0\u05d2\u05e0\u05d5\u05df \u05d6\u05d4 \u05d1\u05ea\u05d5\u05e8 \u05f4\u05d9\u05db\u05d5\u05dc\u05ea \u05e7\u05d9\u05e8\u05d5\u05d1 \u05d3\u05d9\u05e0\u05de\u05d9\u05ea \u05d1\u05e1\u05d2\u05e0\u05d5\u05df UAT\\\" (\u05d4\u05de\u05e9\u05e4\u05d8 \u05d4\u05de\u05e7\u05d5\u05e8\u05d9 \u05de\u05d3\u05d1\u05e8 \u05e2\u05dc \u05d9\u05db\u05d5\u05dc\u05ea \u05e7\u05d9\u05e8\u05d5\u05d1 \u05e1\u05d8\u05d8\u05d9\u05ea \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 ML) \u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d5\u05d3\u05dc \u05de\u05ea\u05d0\u05d9\u05dd \u05ea\u05d5\u05e6\u05d0\u05d4 \u05de\u05ea\u05d0\u05d9\u05de\u05d4 \u05e2\u05dc \u05d1\u05e1\u05d9\u05e1 \u05d4\u05e7\u05dc\u05d8, \u05d5\u05d4\u05ea\u05d5\u05e4\u05e2\u05d4 \u05d4\u05e0\u05e6\u05e4\u05d9\u05ea \u05e0\u05d9\u05ea\u05df \u05dc\u05d4\u05d2\u05d3\u05d9\u05e8 \u05d1\u05ea\u05d5\u05e8 \u05d6\u05d9\u05db\u05e8\u05d5\u05df. \", Word.InsertLocation.replace);\n\n// Paragraph 5: \"Schrodinger's memory\" naming paragraph.\nparagraphs.items[5].insertText(\"\u05d4\u05dd \u05de\u05db\u05e0\u05d9\u05dd \u05d6\u05d0\u05ea \\\"\u05d6\u05d9\u05db\u05e8\u05d5\u05df \u05e9\u05e8\u05d3\u05d9\u05e0\u05d2\u05e8\\\" \u05de\u05db\u05d9\u05d5\u05d5\u05df \u05e9\u05d0\u05e0\u05d5 \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05e7\u05d1\u05d5\u05e2 \u05e9\u05dc-LLMs \u05d9\u05e9 \u05d0\u05ea \u05d4\u05d6\u05d9\u05db\u05e8\u05d5\u05df \u05d4\u05d6\u05d4 \u05e8\u05e7 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05f4\u05e9\u05d0\u05d9\u05dc\u05ea \u05e9\u05d0\u05dc\u05d5\u05ea\u05f4 \u05d5\u05e0\u05d9\u05ea\u05d5\u05d7 \u05d4\u05ea\u05d2\u05d5\u05d1\u05d4 \u05e9\u05dc\u05d5; \u05d0\u05d7\u05e8\u05ea, \u05d4\u05d6\u05d9\u05db\u05e8\u05d5\u05df \u05e0\u05e9\u05d0\u05e8 \u05d1\u05dc\u05ea\u05d9 \u05de\u05d5\u05d2\u05d3\u05e8. \u05d1\u05e0\u05d5\u05e1\u05e3 \u05d1\u05de\u05d0\u05de\u05e8 \u05e0\u05d3\u05d5\u05e0\u05d9\u05dd \u05d2\u05d5\u05e8\u05de\u05d9\u05dd \u05d4\u05de\u05e9\u05e4\u05d9\u05e2\u05d9\u05dd \u05e2\u05dc \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9 LLM: \u05d2\u05d5\u05d3\u05dc \u05d4\u05de\u05d5\u05d3\u05dc, \u05d0\u05d9\u05db\u05d5\u05ea/\u05db\u05de\u05d5\u05ea \u05d4\u05d3\u05d0\u05d8\u05d4 \u05d5\u05d4\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d8\u05d5\u05e2\u05e0\u05d9\u05dd \u05e9\u05d4\u05d6\u05d9\u05db\u05e8\u05d5\u05df \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d1\u05d0\u05d5\u05ea\u05d5 \u05d4\u05d2\u05d5\u05d3\u05dc \u05de\u05d5\u05e9\u05e4\u05e2 \u05de\u05d0\u05d5\u05e4\u05df \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05e9\u05dc\u05d4\u05dd \u05d5\u05d0\u05dd \u05d4\u05de\u05d5\u05d3\u05dc \u05d0\u05d5\u05de\u05df \u05e2\u05dc \u05d9\u05d5\u05ea\u05e8 \u05d3\u05d0\u05d8\u05d4 \u05d0\u05d9\u05db\u05d5\u05ea\u05d9 \u05d0\u05d6 \u05d4\u05d5\u05d0 \u05de\u05e9\u05ea\u05e4\u05e8 (\u05d0\u05d9\u05df \u05d4\u05e4\u05ea\u05e2\u05d5\u05ea \u05db\u05d0\u05df). \", Word.InsertLocation.replace);\n\n// Paragraph 6: previously started with a manual line break; now plain text,\n// no leading break, about the LLM/brain analogy.\nparagraphs.items[6].insertText(\"\u05d5\u05dc\u05d1\u05e1\u05d5\u05e3 \u05e0\u05e2\u05e9\u05d5\u05ea \u05d4\u05e7\u05d1\u05dc\u05d5\u05ea \u05d1\u05d9\u05df \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05ea LLM \u05dc\u05de\u05d1\u05e0\u05d4 \u05d4\u05de\u05d5\u05d3\u05d5\u05dc\u05e8\u05d9 \u05e9\u05dc \u05d4\u05de\u05d5\u05d7 \u05d4\u05d0\u05e0\u05d5\u05e9\u05d9 (\u05d0\u05ea \u05d6\u05d4 \u05e4\u05d7\u05d5\u05ea \u05d0\u05d4\u05d1\u05ea\u05d9 \u05d0\u05d1\u05dc \u05d6\u05e8\u05de\u05ea\u05d9).\", Word.InsertLocation.replace);\n\nawait context.sync();\n\n// Paragraphs 7, 8, 9 (stability/monotonic/\"recommend a look\" paragraphs) are\n// removed entirely in the new version.\nparagraphs.items[9].delete();\nparagraphs.items[8].delete();\nparagraphs.items[7].delete();\nawait context.sync();\n\n// Final paragraph: replace the Nature link with the arXiv link.\nconst remaining = body.paragraphs;\nremaining.load(\"items\");\nawait context.sync();\nconst last = remaining.items[remaining.items.length - 1];\nlast.insertText(\"https://arxiv.org/pdf/2409.10482\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# PowerShell-style Word COM-interop edit script.\n# Applies the same change as edit.js: rewrites paragraphs 1-7 (1-based,\n# i.e. the title block through the \"Schrodinger memory\" naming paragraph and\n# the closing brain-analogy paragraph), removes the three paragraphs that\n# discussed prompt-stability/\"worth a look\" (old paragraphs 8-10), and swaps\n# the trailing Nature link for the arXiv link.\n\n$d = $word.ActiveDocument\n\n# --- Paragraph 1: title line + subtitle line, joined by a manual line break ---\n$d.Paragraphs(1).Range.Text = \"\ud83d\ude80\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 30.09.24: \u26a1\ufe0f\ud83d\ude80\" + [char]11 + \"SCHRODINGER\u2019S MEMORY: LARGE LANGUAGE MODELS\"\n\n# --- Paragraph 2: intro blurb ---\n$d.Paragraphs(2).Range.Text = \"\u05d1\u05d9\u05d5\u05dd \u05d4\u05e1\u05d5\u05e2\u05e8 \u05d4\u05d6\u05d4 (\u05dc\u05de\u05e8\u05d5\u05ea \u05e9\u05d4\u05e1\u05e7\u05d9\u05e8\u05d4 \u05e9\u05d9\u05d9\u05db\u05ea \u05e4\u05d5\u05e8\u05de\u05dc\u05d9\u05ea \u05dc\u05d0\u05ea\u05de\u05d5\u05dc - \u05d0\u05e9\u05dc\u05d9\u05dd \u05d0\u05ea \u05d4\u05e4\u05e2\u05e8 \u05d1\u05d9\u05de\u05d9\u05dd \u05d4\u05e7\u05e8\u05d5\u05d1\u05d9\u05dd) \u05e0\u05e1\u05e7\u05d5\u05e8 \u05de\u05d0\u05de\u05e8 \u05d3\u05d9 \u05e7\u05dc\u05d9\u05dc \u05e2\u05dd \u05e9\u05dd \u05de\u05d0\u05d5\u05d3 \u05dc\u05d0 \u05e7\u05dc\u05d9\u05dc. \u05db\u05d9 \u05d0\u05d9\u05df \u05d3\u05d1\u05e8 \u05e7\u05dc\u05d9\u05dc \u05e9\u05db\u05d5\u05dc\u05dc \u05d1\u05ea\u05d5\u05db\u05d5 \u05d0\u05ea \u05e9\u05de\u05d5 \u05e9\u05dc \u05e9\u05e8\u05d3\u05d9\u05e0\u05d2\u05e8 - \u05e1\u05e4\u05e7 \u05d0\u05dd \u05d4\u05e6\u05dc\u05d7\u05ea\u05d9 \u05dc\u05d4\u05d1\u05d9\u05df \u05d1\u05e6\u05d5\u05e8\u05d4 \u05d8\u05d5\u05d1\u05d4 \u05de\u05e1\u05e4\u05d9\u05e7 \u05d0\u05ea \u05d4\u05de\u05e9\u05d5\u05d5\u05d0\u05d4 \u05e9\u05dc \u05e9\u05e8\u05d3\u05d9\u05e0\u05d2\u05e8 \u05e2\u05d5\u05d3 \u05d1\u05e7\u05d5\u05e8\u05e1 \u05e4\u05d9\u05d6\u05d9\u05e7\u05d4 3 \u05d1\u05d0\u05d5\u05e0\u05d9\u05d1\u05e8\u05e1\u05d9\u05d8\u05d4 \u05d1\u05de\u05d5\u05e1\u05e7\u05d1\u05d4 \u05dc\u05e4\u05e0\u05d9 \u05e2\u05e9\u05e8\u05d5\u05ea \u05e9\u05e0\u05d9\u05dd. \u05d2\u05dd \u05e1\u05d9\u05e4\u05d5\u05e8\u05d5 \u05e9\u05dc \u05d7\u05ea\u05d5\u05dc \u05e9\u05e8\u05d3\u05d9\u05e0\u05d2\u05e8 \u05dc\u05d0 \u05d4\u05ea\u05d1\u05d4\u05e8 \u05e2\u05d3 \u05e2\u05db\u05e9\u05d9\u05d5.\"\n\n# --- Paragraph 3: what the paper covers ---\n$d.Paragraphs(3).Range.Text = \"\u05d0\u05d5\u05e7\u05d9\u05d9, \u05e1\u05d9\u05d9\u05de\u05e0\u05d5 \u05e2\u05dd \u05d4\u05e6\u05d7\u05d5\u05e7\u05d9\u05dd. \u05d4\u05de\u05d0\u05de\u05e8 \u05d7\u05d5\u05e7\u05e8 (\u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05ea) \u05e0\u05d5\u05e9\u05d0 \u05d3\u05d9 \u05e8\u05e6\u05d9\u05e0\u05d9 \u05d5\u05d4\u05d5\u05d0 \u05d4\u05d6\u05db\u05e8\u05d5\u05df \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4. \u05db\u05e9\u05d0\u05e0\u05d7\u05e0\u05d5 \u05e9\u05d5\u05d0\u05dc\u05d9\u05dd LLM \u05de\u05d4 \u05e2\u05d9\u05e8 \u05d4\u05d1\u05d9\u05e8\u05d4 \u05e9\u05dc \u05e9\u05d1\u05d3\u05d9\u05d4, \u05d0\u05d9\u05da \u05d4\u05d5\u05d0 \u05d9\u05d5\u05d3\u05e2 \u05e9\u05d6\u05d4 \u05e1\u05d8\u05d5\u05e7\u05d4\u05d5\u05dc\u05dd. \u05d4\u05de\u05d0\u05de\u05e8 \u05d8\u05d5\u05e2\u05df \u05db\u05d9 \u05d6\u05d9\u05db\u05e8\u05d5\u05df LLM \u05e4\u05d5\u05e2\u05dc \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05ea\u05d0\u05de\u05d4 \u05d3\u05d9\u05e0\u05de\u05d9\u05ea \u05e9\u05dc \u05e4\u05dc\u05d8\u05d9\u05dd \u05dc\u05e7\u05dc\u05d8\u05d9\u05dd. \u05db\u05dc\u05d5\u05de\u05e8 \u05d4\u05de\u05d5\u05d3\u05dc \u05f4\u05d1\u05d5\u05d7\u05e8\u05f4 \u05d0\u05d9\u05da \u05dc\u05e9\u05dc\u05d5\u05e3 \u05d0\u05ea \u05d4\u05de\u05d9\u05d3\u05e2 \u05de\u05d4\u05d6\u05d9\u05db\u05e8\u05d5\u05df \u05d5\u05d1\u05d5\u05e0\u05d4 \u05d0\u05d5\u05ea\u05d5 \u05e2\u05dc \u05e1\u05de\u05da \u05d4\u05e7\u05dc\u05d8. \"\n\n# --- Paragraph 4: transformer / attention explanation ---\n$d.Paragraphs(4).Range.Text = \"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e1\u05d1\u05d9\u05e8\u05d9\u05dd \u05d0\u05ea \u05d0\u05d9\u05da \u05e4\u05d5\u05e2\u05dc \u05d4\u05d6\u05d9\u05db\u05e8\u05d5\u05df \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 \u05e9\u05e4\u05d4 \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea \u05e0\u05d9\u05ea\u05d5\u05d7 \u05e9\u05dc \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05ea \u05d4\u05d8\u05e8\u05e0\u05e1\u05e4\u05d5\u05e8\u05de\u05e8\u05d9\u05dd. \u05de\u05e0\u05d2\u05e0\u05d5\u05df \u05d4-attention (\u05db\u05dc\u05d5\u05de\u05e8 \u05de\u05e7\u05d3\u05de\u05d9 \u05d4-attention \u05e9\u05dc\u05d5) \u05dc\u05de\u05e2\u05e9\u05d4 \u05de\u05d0\u05e4\u05e9\u05e8\u05d9\u05dd \u05dc\u05de\u05d5\u05d3\u05dc \u05dc\u05d1\u05e0\u05d5\u05ea \u05d0\u05ea \u05d4\u05e4\u05dc\u05d8 \u05db\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d4 \u05d3\u05d9\u05e0\u05de\u05d9\u05ea \u05e9\u05dc \u05d4\u05e7\u05dc\u05d8 (\u05db\u05dc\u05d5\u05de\u05e8 \u05dc\u05d0 \u05e7\u05d1\u05d5\u05e2\u05d4 \u05db\u05de\u05d5 \u05d1-MLP \u05d0\u05d5 ConvNets). \"\n\n# --- Paragraph 5: Universal Approximation Theorem paragraph ---\n$d.Paragraphs(5).Range.Text = \"\u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05de\u05e9\u05ea\u05de\u05e9\u05d9\u05dd \u05d1- Universal Approximation Theorem \u05d0\u05d5 UAT \u05db\u05d3\u05d9 \u05dc\u05d4\u05e1\u05d1\u05d9\u05e8 \u05d0\u05ea \u05d4\u05d9\u05db\u05d5\u05dc\u05ea \u05e9\u05dc \u05e9\u05dc\u05d9\u05e4\u05ea \u05de\u05d9\u05d3\u05e2 \u05e9\u05e0\u05dc\u05de\u05d3 \u05d1\u05de\u05d4\u05dc\u05da \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05e2\u05dc \u05d1\u05e1\u05d9\u05e1 \u05ea\u05d5\u05db\u05df \u05e9\u05dc \u05d4\u05e7\u05dc\u05d8. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d8\u05d5\u05e2\u05e0\u05d9\u05dd \u05db\u05d9 \u05e0\u05d9\u05ea\u05df \u05dc\u05d4\u05d1\u05d9\u05df \u05de\u05e0\u05d2\u05e0\u05d5\u05df \u05d6\u05d4 \u05d1\u05ea\u05d5\u05e8 \u05f4\u05d9\u05db\u05d5\u05dc\u05ea \u05e7\u05d9\u05e8\u05d5\u05d1 \u05d3\u05d9\u05e0\u05de\u05d9\u05ea \u05d1\u05e1\u05d2\u05e0\u05d5\u05df UAT`\" (\u05d4\u05de\u05e9\u05e4\u05d8 \u05d4\u05de\u05e7\u05d5\u05e8\u05d9 \u05de\u05d3\u05d1\u05e8 \u05e2\u05dc \u05d9\u05db\u05d5\u05dc\u05ea \u05e7\u05d9\u05e8\u05d5\u05d1 \u05e1\u05d8\u05d8\u05d9\u05ea \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9 ML) \u05db\u05d0\u05e9\u05e8 \u05d4\u05de\u05d5\u05d3\u05dc \u05de\u05ea\u05d0\u05d9\u05dd \u05ea\u05d5\u05e6\u05d0\u05d4 \u05de\u05ea\u05d0\u05d9\u05de\u05d4 \u05e2\u05dc \u05d1\u05e1\u05d9\u05e1 \u05d4\u05e7\u05dc\u05d8, \u05d5\u05d4\u05ea\u05d5\u05e4\u05e2\u05d4 \u05d4\u05e0\u05e6\u05e4\u05d9\u05ea \u05e0\u05d9\u05ea\u05df \u05dc\u05d4\u05d2\u05d3\u05d9\u05e8 \u05d1\u05ea\u05d5\u05e8 \u05d6\u05d9\u05db\u05e8\u05d5\u05df. \"\n\n# --- Paragraph 6: \"Schrodinger's memory\" naming paragraph ---\n$d.Paragraphs(6).Range.Text = \"\u05d4\u05dd \u05de\u05db\u05e0\u05d9\u05dd \u05d6\u05d0\u05ea `\"\u05d6\u05d9\u05db\u05e8\u05d5\u05df \u05e9\u05e8\u05d3\u05d9\u05e0\u05d2\u05e8`\" \u05de\u05db\u05d9\u05d5\u05d5\u05df \u05e9\u05d0\u05e0\u05d5 \u05d9\u05db\u05d5\u05dc\u05d9\u05dd \u05dc\u05e7\u05d1\u05d5\u05e2 \u05e9\u05dc-LLMs \u05d9\u05e9 \u05d0\u05ea \u05d4\u05d6\u05d9\u05db\u05e8\u05d5\u05df \u05d4\u05d6\u05d4 \u05e8\u05e7 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05f4\u05e9\u05d0\u05d9\u05dc\u05ea \u05e9\u05d0\u05dc\u05d5\u05ea\u05f4 \u05d5\u05e0\u05d9\u05ea\u05d5\u05d7 \u05d4\u05ea\u05d2\u05d5\u05d1\u05d4 \u05e9\u05dc\u05d5; \u05d0\u05d7\u05e8\u05ea, \u05d4\u05d6\u05d9\u05db\u05e8\u05d5\u05df \u05e0\u05e9\u05d0\u05e8 \u05d1\u05dc\u05ea\u05d9 \u05de\u05d5\u05d2\u05d3\u05e8. \u05d1\u05e0\u05d5\u05e1\u05e3 \u05d1\u05de\u05d0\u05de\u05e8 \u05e0\u05d3\u05d5\u05e0\u05d9\u05dd \u05d2\u05d5\u05e8\u05de\u05d9\u05dd \u05d4\u05de\u05e9\u05e4\u05d9\u05e2\u05d9\u05dd \u05e2\u05dc \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9 LLM: \u05d2\u05d5\u05d3\u05dc \u05d4\u05de\u05d5\u05d3\u05dc, \u05d0\u05d9\u05db\u05d5\u05ea/\u05db\u05de\u05d5\u05ea \u05d4\u05d3\u05d0\u05d8\u05d4 \u05d5\u05d4\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d4. \u05d4\u05de\u05d7\u05d1\u05e8\u05d9\u05dd \u05d8\u05d5\u05e2\u05e0\u05d9\u05dd \u05e9\u05d4\u05d6\u05d9\u05db\u05e8\u05d5\u05df \u05e9\u05dc \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05d1\u05d0\u05d5\u05ea\u05d5 \u05d4\u05d2\u05d5\u05d3\u05dc \u05de\u05d5\u05e9\u05e4\u05e2 \u05de\u05d0\u05d5\u05e4\u05df \u05d4\u05d0\u05d9\u05de\u05d5\u05df \u05e9\u05dc\u05d4\u05dd \u05d5\u05d0\u05dd \u05d4\u05de\u05d5\u05d3\u05dc \u05d0\u05d5\u05de\u05df \u05e2\u05dc \u05d9\u05d5\u05ea\u05e8 \u05d3\u05d0\u05d8\u05d4 \u05d0\u05d9\u05db\u05d5\u05ea\u05d9 \u05d0\u05d6 \u05d4\u05d5\u05d0 \u05de\u05e9\u05ea\u05e4\u05e8 (\u05d0\u05d9\u05df \u05d4\u05e4\u05ea\u05e2\u05d5\u05ea \u05db\u05d0\u05df). \"\n\n# --- Paragraph 7: previously opened with a manual line break; now plain ---\n# --- text about the LLM / human-brain analogy, no leading break ---\n$d.Paragraphs(7).Range.Text = \"\u05d5\u05dc\u05d1\u05e1\u05d5\u05e3 \u05e0\u05e2\u05e9\u05d5\u05ea \u05d4\u05e7\u05d1\u05dc\u05d5\u05ea \u05d1\u05d9\u05df \u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05ea LLM \u05dc\u05de\u05d1\u05e0\u05d4 \u05d4\u05de\u05d5\u05d3\u05d5\u05dc\u05e8\u05d9 \u05e9\u05dc \u05d4\u05de\u05d5\u05d7 \u05d4\u05d0\u05e0\u05d5\u05e9\u05d9 (\u05d0\u05ea \u05d6\u05d4 \u05e4\u05d7\u05d5\u05ea \u05d0\u05d4\u05d1\u05ea\u05d9 \u05d0\u05d1\u05dc \u05d6\u05e8\u05de\u05ea\u05d9).\"\n\n# --- Remove the three paragraphs on prompt-formulation stability and the ---\n# --- \"worth a look\" recommendation (old paragraphs 8, 9 and 10). Delete from ---\n# --- the highest index down to the lowest, re-resolving Paragraphs(N) fresh ---\n# --- each time, so earlier deletions don't shift the indices out from under us ---\n$d.Paragraphs(10).Range.Delete()\n$d.Paragraphs(9).Range.Delete()\n$d.Paragraphs(8).Range.Delete()\n\n# --- Final paragraph: replace the Nature link with the arXiv link ---\n$d.Paragraphs($d.Paragraphs.Count).Range.Text = \"https://arxiv.org/pdf/2409.10482\"\n\n"}
